$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row1 = @(-7.7485238893292649,-10.058302344363842,-9.5263052770501435,-8.5449785608432194,-9.920824158306301,-9.1988944491296465,-9.3128764409803466,-7.5111190629002209,-9.0453737695867726,-9.2747545587574294,-10.975855091450539,-9.0164502534922732,-9.7598714724370588,-8.7505233023992375,-9.3998620475760752,-9.7261872617058938,-10.839508168100675,-10.150907641910214,-8.4029552227210011,-9.2923434298847312,-9.2260890018487096,-8.6590291824713894,-10.439517746934008,-7.7332097874381542,-9.8373248721556319,-9.375776072466719,-10.038568136980441,-8.1148045724618907,-9.5834841109292768,-8.8901534289369835,-9.9995435496820981,-10.08099416049701,-8.2715780353217436,-9.1887604816239303,-9.1107110459706018,-8.2712542954937742,-9.1950376990014835,-7.3252875667789166,-9.8706672242725944,-10.737726399080485,-9.0635873765650672,-10.626385074754621,-9.2113270394537548,-10.335617418531406,-8.9602472530969965,-8.7420770522866889,-9.842700943203873,-10.034082776379865,-8.910350656420194,-10.070888056171942,-9.2565056468499058,-9.6006926735329152,-8.6551904632578651,-8.7847934781795729,-8.8466299771534516,-9.7833364959127689,-9.8828509695419342,-10.064960424425914,-8.7718411139923287,-9.1039068294630994,-10.406160938129252,-8.9650885076479678,-9.3295223237557519,-8.1482207943170586,-10.782934916803718,-9.8870751336462259,-10.243906551811561,-9.243340797262114,-8.5169651411890293,-10.021694032859189,-9.5841761834975827,-7.8439711082880219,-9.1944009196877712,-8.3231380297072874,-9.749136138772327,-9.6995937031137043,-9.9923344655389084,-9.4842759835411403,-9.077365845149318,-9.2604287991857781,-8.8965214428248167,-10.314259408154101,-9.3207673860724878,-9.8956171785126816,-9.2754888051946907,-10.175457574867924,-9.64477538824978,-8.2249126939601034,-8.6694184968811214,-11.051599776382542,-8.9651210134971322,-10.0959138141742,-9.8827230756280215,-9.4740059239766996,-11.397181003820984,-9.8888803748189193,-10.459304787787126,-9.13082640027568,-9.9922637631551225,-9.7410848830942491)
$row2 = @(-8.7512396780504655,-10.093300817211956,-9.5918708721713202,-8.6437257758507133,-10.045208058728774,-9.240915448424408,-9.3244535705605642,-8.6060544885983425,-9.1057338468491142,-9.2394332287112935,-9.9525012810210143,-9.06206755324237,-9.8533309774453031,-8.766000948187239,-9.4534664446430217,-8.8058706654814056,-10.839965486665919,-10.138108524780472,-9.3507768516489946,-9.2312504496524994,-9.133361955871802,-9.7728055990762392,-10.436732834687852,-8.753508474455959,-9.7541738911933393,-9.4378113612208416,-10.180601736239456,-8.1165311081515377,-9.6495230796199447,-8.8645152576750519,-9.995974961048999,-10.094060120692248,-9.2379720816000432,-8.2501829986555961,-9.0289571211556492,-9.1988284431567475,-9.2736035670958259,-8.3229492056316214,-9.8574434841059979,-9.6547145849981639,-8.029917731635722,-10.596126891147327,-8.1604288514087866,-10.350263616544286,-8.0986610113196829,-9.9086546624840945,-9.7972041805113452,-10.037275623850324,-8.8503351355744169,-9.167489836689704,-9.2728689588426043,-8.5681366453348691,-8.7218519386104685,-8.8026611099458272,-8.8548762251210729,-10.816848161353143,-9.8669173305899314,-9.0597331927181681,-9.7639408290544161,-9.218080981040778,-10.443819208680637,-9.9793062038390428,-9.3392088363322259,-8.1639302376152685,-10.781290801260587,-8.90079494907941,-10.337108694694942,-9.256149932387979,-9.5487375617409924,-10.088973564715179,-9.7031082622422371,-8.8242735935496874,-9.2191023451032148,-9.4288917404963009,-9.7738719254692619,-9.7222831780433587,-10.136380931946467,-9.6105933839906843,-9.0818170086568859,-9.2730980695291656,-10.034172869022926,-10.277498960439873,-9.3455019389072245,-9.8192195920025398,-8.2845526412859378,-10.145895268099338,-9.6461628992699087,-9.2462739182545786,-8.7389005787456107,-11.033454098294166,-9.9934508259496102,-9.1853986339725324,-9.9721081823372675,-9.5639225217822617,-11.288060625384489,-9.8844947444568305,-10.424383423201141,-8.0997408926000478,-9.9360044316868432,-9.8806967511325965)

$data = New-Object "object[,]" 2,100
for ($c = 0; $c -lt 100; $c++) {
    $data[0, $c] = $row1[$c]
    $data[1, $c] = $row2[$c]
}

$ws.Range("A1:CV2").Value = $data
